# Update "浏览/想去人数" (F column) figures across the four sheets to match
# the latest scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 806
$ws.Range("F6").Value = 664
$ws.Range("F7").Value = 1237
$ws.Range("F9").Value = 826
$ws.Range("F10").Value = 701
$ws.Range("F12").Value = 352
$ws.Range("F13").Value = 360
$ws.Range("F14").Value = 724
$ws.Range("F15").Value = 956
$ws.Range("F16").Value = 9998
$ws.Range("F17").Value = 628
$ws.Range("F20").Value = 45
$ws.Range("F22").Value = 274
$ws.Range("F23").Value = 1765
$ws.Range("F26").Value = 489
$ws.Range("F27").Value = 185
$ws.Range("F29").Value = 275
$ws.Range("F30").Value = 193
$ws.Range("F32").Value = 70
$ws.Range("F35").Value = 178
$ws.Range("F36").Value = 197
$ws.Range("F37").Value = 176
$ws.Range("F38").Value = 42

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 135
$ws.Range("F10").Value = 241
$ws.Range("F11").Value = 4440
$ws.Range("F16").Value = 288

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 821

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 821
$ws.Range("F6").Value = 806
$ws.Range("F10").Value = 664
$ws.Range("F11").Value = 1237
$ws.Range("F13").Value = 135
$ws.Range("F14").Value = 826
$ws.Range("F15").Value = 701
$ws.Range("F17").Value = 361
$ws.Range("F19").Value = 956
$ws.Range("F20").Value = 9998
$ws.Range("F21").Value = 241
$ws.Range("F22").Value = 628
$ws.Range("F24").Value = 45
$ws.Range("F25").Value = 274
$ws.Range("F26").Value = 1765
$ws.Range("F28").Value = 489
$ws.Range("F29").Value = 185
$ws.Range("F36").Value = 275
$ws.Range("F37").Value = 193
$ws.Range("F39").Value = 70
$ws.Range("F43").Value = 178
$ws.Range("F46").Value = 197
$ws.Range("F47").Value = 176
